$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Tempo (ms)" column (B) values for the python algorithm chart
$ws.Range("B2").Value  = 0.9279251098632812
$ws.Range("B3").Value  = 1.993179321289062
$ws.Range("B4").Value  = 1.522064208984375
$ws.Range("B5").Value  = 3.943443298339844
$ws.Range("B6").Value  = 2.16984748840332
$ws.Range("B7").Value  = 1.962900161743164
$ws.Range("B8").Value  = 1.932382583618164
$ws.Range("B9").Value  = 1.509904861450195
$ws.Range("B10").Value = 1.974105834960938
$ws.Range("B11").Value = 2.228498458862305
$ws.Range("B12").Value = 2.016425132751465
$ws.Range("B13").Value = 1.968502998352051
